# Regenerate the lattice-multiplication worksheet table with a new set of
# practice problems (commit "Update master to output generated at 503736d").
# Every cell in the 5x3 table keeps its layout (problem / multiplier line /
# rule / two partial-product lines) but gets new numbers, so we just
# rewrite each cell's Range.Text wholesale, using a vertical-tab (chr 11)
# for the manual line breaks Word stores as <w:br/>.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1,1).Range.Text = "13 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"
$t.Cell(1,2).Range.Text = "98 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "9|    |" + $nl + "8|    |"
$t.Cell(1,3).Range.Text = "13 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "1|    |" + $nl + "3|    |"
$t.Cell(2,1).Range.Text = "38 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "3|    |" + $nl + "8|    |"
$t.Cell(2,2).Range.Text = "11 x 56" + $nl + "  5    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "1|    |"
$t.Cell(2,3).Range.Text = "93 x 80" + $nl + "  8    0" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
$t.Cell(3,1).Range.Text = "10 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "0|    |"
$t.Cell(3,2).Range.Text = "42 x 63" + $nl + "  6    3" + $nl + "  ----" + $nl + "4|    |" + $nl + "2|    |"
$t.Cell(3,3).Range.Text = "19 x 96" + $nl + "  9    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "9|    |"
$t.Cell(4,1).Range.Text = "90 x 78" + $nl + "  7    8" + $nl + "  ----" + $nl + "9|    |" + $nl + "0|    |"
$t.Cell(4,2).Range.Text = "21 x 55" + $nl + "  5    5" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"
$t.Cell(4,3).Range.Text = "38 x 71" + $nl + "  7    1" + $nl + "  ----" + $nl + "3|    |" + $nl + "8|    |"
$t.Cell(5,1).Range.Text = "53 x 36" + $nl + "  3    6" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$t.Cell(5,2).Range.Text = "79 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "7|    |" + $nl + "9|    |"
$t.Cell(5,3).Range.Text = "62 x 26" + $nl + "  2    6" + $nl + "  ----" + $nl + "6|    |" + $nl + "2|    |"

Write-Output "All cells updated"